# Re-point the decision-table template from the old "com.myspace" sample
# package to the real project package (com.redhat.prudential_poc), and
# drop the row1/row2 merges so the new (slightly different) values sit in
# a single cell like Excel left them after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Un-merge the header rows that held the old RuleSet/Import values.
$ws.Range("B1:F1").UnMerge()
$ws.Range("B2:F2").UnMerge()

# RuleSet package
$ws.Range("B1").Value = "com.redhat.prudential_poc.rules"

# Import list (Application/Insured POJOs)
$ws.Range("B2").Value = "com.redhat.prudential_poc.pojo.Application,com.redhat.prudential_poc.pojo.Insured"
